$d = $word.ActiveDocument

$d.Content.Find.Execute("920÷7=131, 3", $true, $false, $false, $false, $false, $true, 1, $false, "756÷2=378, 0", 2) | Out-Null
$d.Content.Find.Execute("198÷5=39, 3", $true, $false, $false, $false, $false, $true, 1, $false, "452÷7=64, 4", 2) | Out-Null
$d.Content.Find.Execute("118÷3=39, 1", $true, $false, $false, $false, $false, $true, 1, $false, "251÷4=62, 3", 2) | Out-Null
$d.Content.Find.Execute("779÷3=259, 2", $true, $false, $false, $false, $false, $true, 1, $false, "330÷3=110, 0", 2) | Out-Null
$d.Content.Find.Execute("829÷2=414, 1", $true, $false, $false, $false, $false, $true, 1, $false, "658÷8=82, 2", 2) | Out-Null
$d.Content.Find.Execute("450÷8=56, 2", $true, $false, $false, $false, $false, $true, 1, $false, "341÷9=37, 8", 2) | Out-Null
$d.Content.Find.Execute("362÷4=90, 2", $true, $false, $false, $false, $false, $true, 1, $false, "962÷3=320, 2", 2) | Out-Null
$d.Content.Find.Execute("185÷9=20, 5", $true, $false, $false, $false, $false, $true, 1, $false, "558÷6=93, 0", 2) | Out-Null
$d.Content.Find.Execute("401÷2=200, 1", $true, $false, $false, $false, $false, $true, 1, $false, "976÷5=195, 1", 2) | Out-Null
$d.Content.Find.Execute("321÷6=53, 3", $true, $false, $false, $false, $false, $true, 1, $false, "751÷5=150, 1", 2) | Out-Null
$d.Content.Find.Execute("304÷9=33, 7", $true, $false, $false, $false, $false, $true, 1, $false, "533÷9=59, 2", 2) | Out-Null
$d.Content.Find.Execute("997÷4=249, 1", $true, $false, $false, $false, $false, $true, 1, $false, "840÷6=140, 0", 2) | Out-Null
$d.Content.Find.Execute("686÷3=228, 2", $true, $false, $false, $false, $false, $true, 1, $false, "432÷2=216, 0", 2) | Out-Null
$d.Content.Find.Execute("415÷3=138, 1", $true, $false, $false, $false, $false, $true, 1, $false, "674÷5=134, 4", 2) | Out-Null
$d.Content.Find.Execute("590÷6=98, 2", $true, $false, $false, $false, $false, $true, 1, $false, "387÷9=43, 0", 2) | Out-Null
$d.Content.Find.Execute("668÷4=167, 0", $true, $false, $false, $false, $false, $true, 1, $false, "781÷5=156, 1", 2) | Out-Null
$d.Content.Find.Execute("339÷6=56, 3", $true, $false, $false, $false, $false, $true, 1, $false, "411÷3=137, 0", 2) | Out-Null
$d.Content.Find.Execute("759÷2=379, 1", $true, $false, $false, $false, $false, $true, 1, $false, "552÷8=69, 0", 2) | Out-Null
$d.Content.Find.Execute("302÷8=37, 6", $true, $false, $false, $false, $false, $true, 1, $false, "832÷6=138, 4", 2) | Out-Null
$d.Content.Find.Execute("942÷3=314, 0", $true, $false, $false, $false, $false, $true, 1, $false, "295÷4=73, 3", 2) | Out-Null
$d.Content.Find.Execute("656÷5=131, 1", $true, $false, $false, $false, $false, $true, 1, $false, "772÷3=257, 1", 2) | Out-Null
$d.Content.Find.Execute("856÷2=428, 0", $true, $false, $false, $false, $false, $true, 1, $false, "422÷4=105, 2", 2) | Out-Null
$d.Content.Find.Execute("525÷2=262, 1", $true, $false, $false, $false, $false, $true, 1, $false, "312÷4=78, 0", 2) | Out-Null
$d.Content.Find.Execute("179÷8=22, 3", $true, $false, $false, $false, $false, $true, 1, $false, "618÷2=309, 0", 2) | Out-Null
$d.Content.Find.Execute("171÷8=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "773÷2=386, 1", 2) | Out-Null
